# Convert the M2Doc "m:if / m:else / m:endif" Word field codes used in
# this template into plain-text "{m:...}" runs, mirroring the move to
# TokenIteratorFieldRewriterSplit (which expects literal "{...}" markers
# in the run text rather than legacy fldChar/instrText field codes).

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Escape-XmlText([string]$s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

# Build a run of <w:r><w:t>...</w:t></w:r> elements from an array of text
# pieces - one run per piece, exactly like the template's original
# per-token instrText runs. xml:space="preserve" is only added when the
# piece actually starts/ends with whitespace, matching real Word output.
function Build-RunsXml([string[]]$pieces) {
    $runsXml = ""
    foreach ($piece in $pieces) {
        $escaped = Escape-XmlText $piece
        if ($piece -match '^\s' -or $piece -match '\s$') {
            $runsXml += "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
        } else {
            $runsXml += "<w:r><w:t>$escaped</w:t></w:r>"
        }
    }
    return $runsXml
}

function Wrap-Package([string]$runsXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document ' + $wNs + '><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData>' +
        '</pkg:part></pkg:package>'
}

# The template has exactly 3 legacy field codes, in document order:
#   1. {m:if self.name <> 'anydsl'}
#   2. {m:else}
#   3. {m:endif}
# Replace each field, from the last one to the first one, so that
# deleting/inserting at a later field never shifts the character offsets
# still to be used for the earlier fields.

$replacements = @(
    @("{m:if ", "self.name ", "<>", " ", "'", "anydsl", "'}"),
    @("{m:else}"),
    @("{m:endif}")
)

$fieldCount = $d.Fields.Count
for ($i = $fieldCount; $i -ge 1; $i--) {
    $f = $d.Fields($i)
    $pos = $f.Code.Start - 1
    $f.Delete()

    $runsXml = Build-RunsXml $replacements[$i - 1]
    $pkgXml = Wrap-Package $runsXml
    $r = $d.Range($pos, $pos)
    $r.InsertXML($pkgXml) | Out-Null
}
